$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"0.020335"
$ws.Range("H2").Value = [double]"0.061005"
$ws.Range("I2").Value = [double]"0.009804808687698561"
$ws.Range("J2").Value = [double]"0.009804808687698559"
$ws.Range("M2").Value = [double]"1.306600666666667"
$ws.Range("N2").Value = [double]"3.919802"
$ws.Range("O2").Value = [double]"0.8137131711319011"
$ws.Range("P2").Value = [double]"0.8137131711319011"
$ws.Range("Q2").Value = [double]"0.02656972455666666"
$ws.Range("R2").Value = [double]"0.23912752101"
$ws.Range("S2").Value = [double]"0.00797830196960881"
$ws.Range("T2").Value = [double]"0.007978301969608808"
$ws.Range("G3").Value = [double]"0.020335"
$ws.Range("H3").Value = [double]"0.061005"
$ws.Range("I3").Value = [double]"0.009804808687698561"
$ws.Range("J3").Value = [double]"0.009804808687698559"
$ws.Range("O3").Value = [double]"0.006864390964089149"
$ws.Range("P3").Value = [double]"0.006864390964089149"
$ws.Range("Q3").Value = [double]"0.0002241391483333333"
$ws.Range("R3").Value = [double]"0.002017252335"
$ws.Range("S3").Value = [double]"6.730404016046079E-05"
$ws.Range("T3").Value = [double]"6.730404016046078E-05"
$ws.Range("G4").Value = [double]"0.020335"
$ws.Range("H4").Value = [double]"0.061005"
$ws.Range("I4").Value = [double]"0.009804808687698561"
$ws.Range("J4").Value = [double]"0.009804808687698559"
$ws.Range("M4").Value = [double]"0.2881033333333333"
$ws.Range("N4").Value = [double]"0.86431"
$ws.Range("O4").Value = [double]"0.1794224379040098"
$ws.Range("P4").Value = [double]"0.1794224379040098"
$ws.Range("Q4").Value = [double]"0.005858581283333333"
$ws.Range("R4").Value = [double]"0.05272723155"
$ws.Range("S4").Value = [double]"0.001759202677929291"
$ws.Range("T4").Value = [double]"0.001759202677929291"
$ws.Range("I5").Value = [double]"0.1486140913768632"
$ws.Range("J5").Value = [double]"0.1486140913768632"
$ws.Range("M5").Value = [double]"1.306600666666667"
$ws.Range("N5").Value = [double]"3.919802"
$ws.Range("O5").Value = [double]"0.8137131711319011"
$ws.Range("P5").Value = [double]"0.8137131711319011"
$ws.Range("Q5").Value = [double]"0.402724377282"
$ws.Range("R5").Value = [double]"3.624519395538"
$ws.Range("S5").Value = [double]"0.1209292435691535"
$ws.Range("T5").Value = [double]"0.1209292435691535"
$ws.Range("I6").Value = [double]"0.1486140913768632"
$ws.Range("J6").Value = [double]"0.1486140913768632"
$ws.Range("O6").Value = [double]"0.006864390964089149"
$ws.Range("P6").Value = [double]"0.006864390964089149"
$ws.Range("S6").Value = [double]"0.001020145225983659"
$ws.Range("T6").Value = [double]"0.001020145225983659"
$ws.Range("I7").Value = [double]"0.1486140913768632"
$ws.Range("J7").Value = [double]"0.1486140913768632"
$ws.Range("M7").Value = [double]"0.2881033333333333"
$ws.Range("N7").Value = [double]"0.86431"
$ws.Range("O7").Value = [double]"0.1794224379040098"
$ws.Range("P7").Value = [double]"0.1794224379040098"
$ws.Range("Q7").Value = [double]"0.08880007371000001"
$ws.Range("R7").Value = [double]"0.7992006633900001"
$ws.Range("S7").Value = [double]"0.02666470258172608"
$ws.Range("T7").Value = [double]"0.02666470258172608"
$ws.Range("H8").Value = [double]"5.236273000000001"
$ws.Range("I8").Value = [double]"0.8415810999354383"
$ws.Range("J8").Value = [double]"0.8415810999354382"
$ws.Range("M8").Value = [double]"1.306600666666667"
$ws.Range("N8").Value = [double]"3.919802"
$ws.Range("O8").Value = [double]"0.8137131711319011"
$ws.Range("P8").Value = [double]"0.8137131711319011"
$ws.Range("Q8").Value = [double]"2.280572597549555"
$ws.Range("R8").Value = [double]"20.525153377946"
$ws.Range("S8").Value = [double]"0.6848056255931388"
$ws.Range("T8").Value = [double]"0.6848056255931387"
$ws.Range("H9").Value = [double]"5.236273000000001"
$ws.Range("I9").Value = [double]"0.8415810999354383"
$ws.Range("J9").Value = [double]"0.8415810999354382"
$ws.Range("O9").Value = [double]"0.006864390964089149"
$ws.Range("P9").Value = [double]"0.006864390964089149"
$ws.Range("S9").Value = [double]"0.00577694169794503"
$ws.Range("T9").Value = [double]"0.005776941697945029"
$ws.Range("H10").Value = [double]"5.236273000000001"
$ws.Range("I10").Value = [double]"0.8415810999354383"
$ws.Range("J10").Value = [double]"0.8415810999354382"
$ws.Range("M10").Value = [double]"0.2881033333333333"
$ws.Range("N10").Value = [double]"0.86431"
$ws.Range("O10").Value = [double]"0.1794224379040098"
$ws.Range("P10").Value = [double]"0.1794224379040098"
$ws.Range("S10").Value = [double]"0.1509985326443544"
$ws.Range("T10").Value = [double]"0.1509985326443544"
